# Auto-generated edit script: updates crypto price/volume table cells
# per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.037.17"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "3.464.22"
$ws.Range("E3").Value = "  +4.28%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'585.43"
$ws.Range("E5").Value = "  +6.09%  "
$ws.Range("D6").Value = "'187.27"
$ws.Range("E6").Value = "  +8.48%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("D8").Value = "3.457.50"
$ws.Range("E8").Value = "  +4.43%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'0.648"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("D12").Value = "'56.18"
$ws.Range("E12").Value = "  +5.97%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'9.43"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("D15").Value = "4.021.93"
$ws.Range("E15").Value = "  +4.38%  "
$ws.Range("D16").Value = "'18.79"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").Value = "3.462.48"
$ws.Range("E17").Value = "  +4.43%  "
$ws.Range("D18").Value = "66.972.53"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("D19").Value = "'12.17"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'486.75"
$ws.Range("E22").Value = "  +7.98%  "
$ws.Range("D23").Value = "'5.26"
$ws.Range("E23").Value = "  +5.28%  "
$ws.Range("D24").Value = "'16.74"
$ws.Range("E24").Value = "  +20.91%  "
$ws.Range("D25").Value = "'4.50"
$ws.Range("E25").Value = "  +11.40%  "
$ws.Range("D26").Value = "'89.65"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").Value = "'2.94"
$ws.Range("E27").Value = "  +3.24%  "
$ws.Range("D28").Value = "'10.96"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "'9.13"
$ws.Range("E29").Value = "  +6.75%  "
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("E31").Value = "  +10.96%  "
$ws.Range("D32").Value = "'599.09"
$ws.Range("E32").Value = "  +5.28%  "
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").Value = "'63.86"
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("D36").Value = "'0.149"
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D38").Value = "'36.49"
$ws.Range("E38").Value = "  +3.91%  "
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "'0.385"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0757"
$ws.Range("E41").Value = "  +4.41%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.246.69"
$ws.Range("E42").Value = "  +6.25%  "
$ws.Range("E43").Value = "  +6.68%  "
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.29"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.52"
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  +23.55%  "
$ws.Range("D48").Value = "'0.136"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").Value = "'3.27"
$ws.Range("E49").Value = "  +13.69%  "
$ws.Range("E50").Value = "  +7.07%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.07%  "
